# Word COM-interop edit script
# 1) Widen the two dashed separator lines in the R console output block (cosmetic re-wrap)
# 2) Rewrite the OTU-threshold paragraph: drop the placeholder sentence, fill in the
#    actual threshold values (bolded) for the full-length/V4/V3-V4/V4-V5 regions.

$d = $word.ActiveDocument

$sep1 = $d.Content
$okSep1 = $sep1.Find.Execute("## ── Attaching packages ─────────────────────────────────────────────────────────── tidyverse 1.3.0 ──", $false, $false, $false, $false, $false, $true, 1, $false, "## ── Attaching packages ──────────────────────────────────────────────────────────────────────────────────────── tidyverse 1.3.0 ──", 2)
Write-Output "Separator 1 replaced: $okSep1"

$sep2 = $d.Content
$okSep2 = $sep2.Find.Execute("## ── Conflicts ────────────────────────────────────────────────────────────── tidyverse_conflicts() ──", $false, $false, $false, $false, $false, $true, 1, $false, "## ── Conflicts ─────────────────────────────────────────────────────────────────────────────────────────── tidyverse_conflicts() ──", 2)
Write-Output "Separator 2 replaced: $okSep2"

# Locate the OTU-threshold paragraph by its (pre-edit) exact text
$rng = $d.Content
$okPara = $rng.Find.Execute("A method to avoid splitting a single genome into multiple units of inference is to cluster 16S rRNA gene sequences together that are similar. However, this also increases the risk of lumping together genes from different species that are similar to each other. Therefore, I assessed the impact of the threshold used to define clusters of 16S rRNA genes on the propopensity to lump species together and split genome apart. I identified the threshold where 90% of bacterial species would be represented by a single OTU. For full length 16S rRNA gene sequences, I found that at a threshold of XX%, 90% of the species would be represented by a single OTU. Similarly, thresholds of XX, XX, and XX% were observed for the V4, V3-V4, and V4-V5 regions. However, at these thresholds, multiple species could be represented by the same OTU. At the highest level of resolution, XX% of the species shared a 16S rRNA gene sequence variant with another species. Given the risk of splitting a genome into multiple OTUs is more biologically problematic than lumping species together, larger thresholds are advisable.")
Write-Output "Paragraph located: $okPara"

# Clear the old sentence, then splice in the new runs (plain / bold-number / plain, …)
# via InsertXML so each segment keeps its own <w:r> even where formatting repeats,
# matching how the source document structures bolded inline values elsewhere.
$rng.Text = ""
[void]$rng.InsertXML("<pkg:package xmlns:pkg=""http://schemas.microsoft.com/office/2006/xmlPackage""><pkg:part pkg:name=""/word/document.xml"" pkg:contentType=""application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml""><pkg:xmlData><w:document xmlns:w=""http://schemas.openxmlformats.org/wordprocessingml/2006/main""><w:body><w:p><w:pPr><w:pStyle w:val=""BodyText""/></w:pPr><w:r><w:t xml:space=""preserve"">A method to avoid splitting a single genome into multiple units of inference is to cluster 16S rRNA gene sequences together that are similar. However, this also increases the risk of lumping together genes from different species that are similar to each other. Therefore, I assessed the impact of the threshold used to define clusters of 16S rRNA genes on the propensity to split a genome apart or to lump species together. For full length 16S rRNA gene sequences, I found that at a threshold of</w:t></w:r><w:r><w:t xml:space=""preserve""> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=""preserve"">5.5</w:t></w:r><w:r><w:t xml:space=""preserve"">%, 95% of the species with 7 copies of the rrn operon would be represented by a single OTU. Similarly, thresholds of</w:t></w:r><w:r><w:t xml:space=""preserve""> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=""preserve"">2.5</w:t></w:r><w:r><w:t xml:space=""preserve"">,</w:t></w:r><w:r><w:t xml:space=""preserve""> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=""preserve"">4.0</w:t></w:r><w:r><w:t xml:space=""preserve"">, and</w:t></w:r><w:r><w:t xml:space=""preserve""> </w:t></w:r><w:r><w:rPr><w:b/></w:rPr><w:t xml:space=""preserve"">3.5</w:t></w:r><w:r><w:t xml:space=""preserve"">% were observed for the V4, V3-V4, and V4-V5 regions, respectively. However, at these thresholds, multiple species could be represented by the same OTU. At the highest level of resolution, XX% of the species shared a 16S rRNA gene sequence variant with another species. Given the risk of splitting a genome into multiple OTUs is more biologically problematic than lumping species together, larger thresholds are advisable.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>")
Write-Output "Paragraph rewritten."
